$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "64.711.11"; E = "  +0.81%  " }
    @{ Row = 3; D = "3.153.89"; E = "  +2.16%  " }
    @{ Row = 4; D = $null; E = "  -0.09%  " }
    @{ Row = 5; D = "570.44"; E = "  +2.00%  " }
    @{ Row = 6; D = "150.70"; E = "  +4.68%  " }
    @{ Row = 7; D = $null; E = "  -0.05%  " }
    @{ Row = 8; D = "3.153.67"; E = "  +2.28%  " }
    @{ Row = 9; D = $null; E = "  +4.24%  " }
    @{ Row = 10; D = $null; E = "  +5.02%  " }
    @{ Row = 11; D = "6.15"; E = "  +0.24%  " }
    @{ Row = 12; D = "0.504"; E = "  +6.80%  " }
    @{ Row = 13; D = $null; E = "  +13.23%  " }
    @{ Row = 14; D = $null; E = "  +8.78%  " }
    @{ Row = 15; D = "3.667.83"; E = "  +1.80%  " }
    @{ Row = 16; D = "64.829.27"; E = "  +0.81%  " }
    @{ Row = 17; D = $null; E = "  +7.12%  " }
    @{ Row = 18; D = "3.152.52"; E = "  +1.79%  " }
    @{ Row = 19; D = $null; E = "  +0.38%  " }
    @{ Row = 20; D = "515.84"; E = "  +7.18%  " }
    @{ Row = 21; D = "14.97"; E = "  +6.66%  " }
    @{ Row = 22; D = $null; E = "  +8.68%  " }
    @{ Row = 23; D = "15.32"; E = "  +8.53%  " }
    @{ Row = 24; D = $null; E = "  +4.22%  " }
    @{ Row = 25; D = "85.19"; E = "  +4.69%  " }
    @{ Row = 26; D = $null; E = "  +0.09%  " }
    @{ Row = 27; D = $null; E = "  +4.23%  " }
    @{ Row = 28; D = "8.80"; E = "  +10.00%  " }
    @{ Row = 29; D = $null; E = "  +6.75%  " }
    @{ Row = 30; D = "27.93"; E = "  +6.05%  " }
    @{ Row = 31; D = "1.00"; E = "  -0.05%  " }
    @{ Row = 32; D = $null; E = "  +4.65%  " }
    @{ Row = 33; D = "2.67"; E = "  +8.04%  " }
    @{ Row = 34; D = "6.17"; E = "  +9.30%  " }
    @{ Row = 35; D = $null; E = "  +6.36%  " }
    @{ Row = 36; D = "55.71"; E = "  +0.41%  " }
    @{ Row = 37; D = "487.11"; E = "  +10.19%  " }
    @{ Row = 38; D = $null; E = "  +6.64%  " }
    @{ Row = 39; D = $null; E = "  +3.96%  " }
    @{ Row = 40; D = "2.97"; E = "  +0.09%  " }
    @{ Row = 41; D = "3.117.84"; E = "  +5.46%  " }
    @{ Row = 42; D = $null; E = "  +5.65%  " }
    @{ Row = 43; D = $null; E = "  +5.02%  " }
    @{ Row = 44; D = $null; E = "  +12.64%  " }
    @{ Row = 45; D = $null; E = "  +15.44%  " }
    @{ Row = 46; D = "29.50"; E = "  +4.76%  " }
    @{ Row = 47; D = $null; E = "  +12.00%  " }
    @{ Row = 49; D = $null; E = "  +3.10%  " }
    @{ Row = 50; D = $null; E = "  +10.50%  " }
    @{ Row = 51; D = "119.64"; E = "  +1.08%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
